# Weekly data refresh: insert a new price-report row for Mango at
# "Feria Lagunitas de Puerto Montt" above the current first data row (row 381).
# All existing data rows shift down by one (381->382, ..., 428->429).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 381, pushing rows 381..428 down to 382..429.
$ws.Rows.Item(381).Insert()

# Populate the new row 381 with this week's report.
$ws.Range("A381").Value = 4
$ws.Range("B381").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C381").Value = "Los Lagos"
$ws.Range("D381").Value = 45154
$ws.Range("E381").Value = 10
$ws.Range("F381").Value = "Fruta"
$ws.Range("G381").Value = 100108
$ws.Range("H381").Value = "Tropicales y subtropicales"
$ws.Range("I381").Value = 100108002
$ws.Range("J381").Value = "Mango"
$ws.Range("K381").Value = "Sin especificar"
$ws.Range("L381").Value = "Primera"
$ws.Range("M381").Value = 120
$ws.Range("N381").Value = 10000
$ws.Range("O381").Value = 11000
$ws.Range("P381").Value = 10500
$ws.Range("Q381").Value = "$/bandeja 4 kilos"
$ws.Range("R381").Value = "Perú"
$ws.Range("S381").Value = 2625
$ws.Range("T381").Value = 4
